# Gacha.xlsx - "중복 시 획득 마일리지 시트 추가(241109)"
# Adds a new "When_Dup" worksheet (Item_Grade / Acquired_Mileage) after the
# existing sheets, with a header row styled like the other sheets, a
# threaded note on A1 explaining the grade codes, and makes it the active
# sheet/selection the way the authored workbook ends up.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "When_Dup" sheet as the last tab.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "When_Dup"

# ---------------------------------------------------------------------
# 2. Header + data.
# ---------------------------------------------------------------------
$newSheet.Range("A1").Value = "Item_Grade"
$newSheet.Range("B1").Value = "Acquired_Mileage"

$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 10
$newSheet.Range("A3").Value = 2
$newSheet.Range("B3").Value = 20
$newSheet.Range("A4").Value = 3
$newSheet.Range("B4").Value = 30

# ---------------------------------------------------------------------
# 3. Match header styling (bold font / yellow fill / centered) used by
#    the other sheets' header rows, by copying the format from GachaGroup!A1:B1.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:B1").Copy()
$newSheet.Range("A1:B1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Column widths (12 / 18 characters, matching the authored sheet).
# ---------------------------------------------------------------------
$newSheet.Columns.Item(1).ColumnWidth = 11.285714285714286
$newSheet.Columns.Item(2).ColumnWidth = 17.285714285714285

# ---------------------------------------------------------------------
# 5. Threaded comment on A1 documenting the grade codes.
# ---------------------------------------------------------------------
$null = $newSheet.Range("A1").AddCommentThreaded("Normal = 1`nRare = 2`nUnique = 3")

# ---------------------------------------------------------------------
# 6. Final selection/activation state: GachaGroup_Item keeps a selection
#    at J18 (no longer the active tab), and the new When_Dup sheet becomes
#    the active tab with its selection at F21.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("GachaGroup_Item")
$ws2.Activate()
$ws2.Range("J18").Select()

$newSheet.Activate()
$newSheet.Range("F21").Select()
